$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Motley"
$ws.Range("C3").Value = "Lorrie"
$ws.Range("G3").Value = "15 Lewis Ave apt b2"
$ws.Range("I3").Value = "Salem, VA 24153"

$ws.Range("J6").Select()
